$d = $word.ActiveDocument

# The "Bibliografia" paragraph originally contains a single run whose text
# concatenates 8 bibliography entries with no separators. Insert a manual
# line break (^l -> <w:br/>) between each entry, splitting the run into
# multiple runs joined by <w:br/> elements, matching the target diff.

$d.Content.Find.Execute(
    "Cengage, 2ed, 2022.Nilo Ney Coutinho Menezes.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Cengage, 2ed, 2022.^lNilo Ney Coutinho Menezes.", 2)

$d.Content.Find.Execute(
    "Para Iniciantes, 3a ed, 2019.Ramalho, L.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Para Iniciantes, 3a ed, 2019.^lRamalho, L.", 2)

$d.Content.Find.Execute(
    "O’Reilly-Novatec, 2015Downey, A. B.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "O’Reilly-Novatec, 2015^lDowney, A. B.", 2)

$d.Content.Find.Execute(
    "O’Reilly-Novatec, 2016.STEWART, J. M.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "O’Reilly-Novatec, 2016.^lSTEWART, J. M.", 2)

$d.Content.Find.Execute(
    "Cambridge University Press, 2014.TELLES, M.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Cambridge University Press, 2014.^lTELLES, M.", 2)

$d.Content.Find.Execute(
    "Thomson Course Technology PTR, 2008.LUTZ, Mark.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Thomson Course Technology PTR, 2008.^lLUTZ, Mark.", 2)

$d.Content.Find.Execute(
    "O’Reilly Media, 2006.MCGREGGOR, D. M.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "O’Reilly Media, 2006.^lMCGREGGOR, D. M.", 2)
